$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: shift the existing quarters down
#    one row and insert the new "2022-Q4" figures at the top of the
#    data block.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Range("B7").Value = "2020-Q4"
$summary.Range("C7").Value = 14
$summary.Range("D7").Value = 19.95

$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 4.41

$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 19
$summary.Range("D5").Value = 3.19

$summary.Range("B4").Value = "2021-Q3"
$summary.Range("C4").Value = 7
$summary.Range("D4").Value = 1.1

$summary.Range("B3").Value = "2021-Q4"
$summary.Range("C3").Value = 9
$summary.Range("D3").Value = 5.01

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 1.26

# Give the newly created A7 cell the same formatting as the rest of
# column A (index marker column), then set its value.
$summary.Range("A2").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q4" worksheet right before the existing
#    "2021-Q4" sheet, holding the quarterly fund holdings table.
#    We duplicate the "2021-Q4" sheet (so all formatting/styles carry
#    over cleanly) and then overwrite its contents, trimming the extra
#    rows it doesn't need.
# ---------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)    # "2021-Q4"
$refSheet.Copy($refSheet)             # duplicate placed right before it
$newSheet = $wb.Worksheets.Item(2)    # the freshly inserted duplicate
$newSheet.Name = "2022-Q4"

# The template sheet has 10 data rows (2..11); the new table only needs
# 5 (2..6), so drop the extra ones entirely.
$newSheet.Range("A7:H11").Delete()

# Header row text.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$data = @(
    @(0, "090007", "大成策略回报混合", "14.14", "60.14", "3.74", "0.5288", 4),
    @(1, "008269", "大成睿享混合A", "23.67", "63.04", "2.14", "0.5065", 10),
    @(2, "008270", "大成睿享混合C", "8.29", "63.04", "2.14", "0.1774", 10),
    @(3, "015564", "大成弘远回报一年持有混合A", "2.54", "27.63", "1.91", "0.0485", 6),
    @(4, "015565", "大成弘远回报一年持有混合C", "0.09", "27.63", "1.91", "0.0017", 6)
)

$row = 2
foreach ($entry in $data) {
    $newSheet.Cells.Item($row, 1).Value = $entry[0]
    # Fund code / size / position columns store numeric-looking values
    # as plain text (matching the rest of the workbook), so force text
    # via the NumberFormat before assignment.
    $newSheet.Cells.Item($row, 2).NumberFormat = "@"
    $newSheet.Cells.Item($row, 2).Value = $entry[1]
    $newSheet.Cells.Item($row, 3).Value = $entry[2]
    $newSheet.Cells.Item($row, 4).NumberFormat = "@"
    $newSheet.Cells.Item($row, 4).Value = $entry[3]
    $newSheet.Cells.Item($row, 5).NumberFormat = "@"
    $newSheet.Cells.Item($row, 5).Value = $entry[4]
    $newSheet.Cells.Item($row, 6).NumberFormat = "@"
    $newSheet.Cells.Item($row, 6).Value = $entry[5]
    $newSheet.Cells.Item($row, 7).NumberFormat = "@"
    $newSheet.Cells.Item($row, 7).Value = $entry[6]
    $newSheet.Cells.Item($row, 8).Value = $entry[7]
    $row = $row + 1
}

Write-Output "done"
